$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep text formatting (values look numeric but are stored as text,
# matching the original inlineStr cell type in the workbook).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.836.40'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.78%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.903.64'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.96%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '568.72'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -4.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.35'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.40%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.510'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.39%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.902.74'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.91%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -8.64%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.59%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.436'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.02%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000232'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -3.12%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.19'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.85%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.82%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.385.22'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.93%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.795.99'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.68%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.03%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.906.30'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.71%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '437.25'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.45%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.33'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.659'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.55%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.94'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.53%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '79.60'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -2.44%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.86'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.17'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -9.32%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.04'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -4.44%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000104'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +7.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.11'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.86%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.51'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -4.51%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.07'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -4.29%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.107'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.74%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.23%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '25.71'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -3.11%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.51%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '49.01'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.11%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -4.52%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.82'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -9.57%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.27%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.30'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.67%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.25%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.270'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -4.16%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.702.19'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.56%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '133.94'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.35%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '338.88'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -6.66%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.104'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.09%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '21.82'
